$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 6.274549265160829
$ws.Range("C2").Value = 1.119709817049454
$ws.Range("D2").Value = 0.3414738407473408
$ws.Range("E2").Value = 0.0581337311815362
$ws.Range("G2").Value = 0.0027381969358049
$ws.Range("I2").Value = 3.906471369855325
$ws.Range("J2").Value = 0.01233293781549349
$ws.Range("L2").Value = 1.035067476419357
$ws.Range("N2").Value = 3.519931361770077
$ws.Range("B3").Value = 6.154107944431189
$ws.Range("C3").Value = 1.078568762718874
$ws.Range("D3").Value = 0.3415760845195379
$ws.Range("E3").Value = 0.05832720957398907
$ws.Range("G3").Value = 0.002746303904002672
$ws.Range("I3").Value = 3.886125821971802
$ws.Range("J3").Value = 0.0108928029096127
$ws.Range("L3").Value = 1.026592235534807
$ws.Range("N3").Value = 3.525171018865962
$ws.Range("B4").Value = 6.084097407814625
$ws.Range("C4").Value = 1.054079031881997
$ws.Range("D4").Value = 0.3417794674237555
$ws.Range("E4").Value = 0.05845424535653798
$ws.Range("G4").Value = 0.002751537759493231
$ws.Range("I4").Value = 3.875349539574998
$ws.Range("J4").Value = 0.01000586258249925
$ws.Range("L4").Value = 1.021939671701858
$ws.Range("N4").Value = 3.529199036031372
$ws.Range("B5").Value = 6.056553422898446
$ws.Range("C5").Value = 1.044291145164493
$ws.Range("D5").Value = 0.3418976752420946
$ws.Range("E5").Value = 0.05850809059227569
$ws.Range("G5").Value = 0.002753735262990054
$ws.Range("I5").Value = 3.871387498263431
$ws.Range("J5").Value = 0.009643684793143592
$ws.Range("L5").Value = 1.02018183778398
$ws.Range("N5").Value = 3.531043292702051
$ws.Range("B6").Value = 6.05203916328179
$ws.Range("C6").Value = 1.042677405004213
$ws.Range("D6").Value = 0.341919436479273
$ws.Range("E6").Value = 0.05851715716240502
$ws.Range("G6").Value = 0.002754104069857541
$ws.Range("I6").Value = 3.87075547690732
$ws.Range("J6").Value = 0.0095834987063661
$ws.Range("L6").Value = 1.019898278885663
$ws.Range("N6").Value = 3.531361747594971
$ws.Range("B7").Value = 6.083721954836051
$ws.Range("C7").Value = 1.053946254735251
$ws.Range("D7").Value = 0.3417809186101124
$ws.Range("E7").Value = 0.05845496311486831
$ws.Range("G7").Value = 0.002751567133580028
$ws.Range("I7").Value = 3.875294370299699
$ws.Range("J7").Value = 0.01000098121834725
$ws.Range("L7").Value = 1.0219154063236
$ws.Range("N7").Value = 3.52922308858362
$ws.Range("B8").Value = 6.232199997698899
$ws.Range("C8").Value = 1.105363117636841
$ws.Range("D8").Value = 0.3414798932076053
$ws.Range("E8").Value = 0.05819873564103117
$ws.Range("G8").Value = 0.002740939209826375
$ws.Range("I8").Value = 3.899098778842216
$ws.Range("J8").Value = 0.01183689790951803
$ws.Range("L8").Value = 1.032030488877723
$ws.Range("N8").Value = 3.521569044022527
$ws.Range("B9").Value = 6.55489330453139
$ws.Range("C9").Value = 1.212407644043765
$ws.Range("D9").Value = 0.3420069655094693
$ws.Range("E9").Value = 0.05776141464833762
$ws.Range("G9").Value = 0.002722118593607158
$ws.Range("I9").Value = 3.959497614228198
$ws.Range("J9").Value = 0.01541889847289113
$ws.Range("L9").Value = 1.0562655446181
$ws.Range("N9").Value = 3.51304368885657
$ws.Range("B10").Value = 6.811597114005167
$ws.Range("C10").Value = 1.294990338404489
$ws.Range("D10").Value = 0.3430784786727088
$ws.Range("E10").Value = 0.05747950781332234
$ws.Range("G10").Value = 0.002709506570086607
$ws.Range("I10").Value = 4.012389805908256
$ws.Range("J10").Value = 0.0180440703575897
$ws.Range("L10").Value = 1.076792384527039
$ws.Range("N10").Value = 3.510805158198195
$ws.Range("B11").Value = 6.932730474160053
$ws.Range("C11").Value = 1.333447487769433
$ws.Range("D11").Value = 0.3437152814320541
$ws.Range("E11").Value = 0.05735974850729253
$ws.Range("G11").Value = 0.002704029442509675
$ws.Range("I11").Value = 4.038336054896263
$ws.Range("J11").Value = 0.01923793211065572
$ws.Range("L11").Value = 1.086730537072896
$ws.Range("N11").Value = 3.510676810138392
$ws.Range("B12").Value = 6.979234164871514
$ws.Range("C12").Value = 1.348140760433068
$ws.Range("D12").Value = 0.343977958964544
$ws.Range("E12").Value = 0.05731561332024659
$ws.Range("G12").Value = 0.002701992532581925
$ws.Range("I12").Value = 4.048435085518832
$ws.Range("J12").Value = 0.0196900471272059
$ws.Range("L12").Value = 1.090580893312136
$ws.Range("N12").Value = 3.510757506493235
$ws.Range("B13").Value = 6.969190499159367
$ws.Range("C13").Value = 1.344970463164543
$ws.Range("D13").Value = 0.3439204280007147
$ws.Range("E13").Value = 0.05732506465432685
$ws.Range("G13").Value = 0.002702429568838995
$ws.Range("I13").Value = 4.046247860706615
$ws.Range("J13").Value = 0.01959267373235463
$ws.Range("L13").Value = 1.089747770429767
$ws.Range("N13").Value = 3.510734355132939
$ws.Range("B14").Value = 6.93654363418716
$ws.Range("C14").Value = 1.334653684957743
$ws.Range("D14").Value = 0.3437364601637825
$ws.Range("E14").Value = 0.05735609315304924
$ws.Range("G14").Value = 0.002703861121298052
$ws.Range("I14").Value = 4.039161406516627
$ws.Range("J14").Value = 0.01927512694160072
$ws.Range("L14").Value = 1.0870455608692
$ws.Range("N14").Value = 3.510680850072163
$ws.Range("B15").Value = 6.916629133798892
$ws.Range("C15").Value = 1.328351413289283
$ws.Range("D15").Value = 0.343626580640219
$ws.Range("E15").Value = 0.05737525710753921
$ws.Range("G15").Value = 0.002704742820408163
$ws.Range("I15").Value = 4.034856478468967
$ws.Range("J15").Value = 0.01908062592163162
$ws.Range("L15").Value = 1.085401728525341
$ws.Range("N15").Value = 3.510664954886749
$ws.Range("B16").Value = 6.803769000367424
$ws.Range("C16").Value = 1.292495185449411
$ws.Range("D16").Value = 0.343039872405086
$ws.Range("E16").Value = 0.05748750464301766
$ws.Range("G16").Value = 0.002709869726844267
$ws.Range("I16").Value = 4.010732319391948
$ws.Range("J16").Value = 0.01796604815295666
$ws.Range("L16").Value = 1.076155040039367
$ws.Range("N16").Value = 3.510831566982375
$ws.Range("B17").Value = 6.735653597216128
$ws.Range("C17").Value = 1.270728173061343
$ws.Range("D17").Value = 0.3427182405294644
$ws.Range("E17").Value = 0.05755853391579313
$ws.Range("G17").Value = 0.002713081369909004
$ws.Range("I17").Value = 3.996417664958486
$ws.Range("J17").Value = 0.01728225588049526
$ws.Range("L17").Value = 1.070636731937071
$ws.Range("N17").Value = 3.511162653755889
$ws.Range("B18").Value = 6.696885166978632
$ws.Range("C18").Value = 1.258291973372536
$ws.Range("D18").Value = 0.3425473018947685
$ws.Range("E18").Value = 0.05760018666468802
$ws.Range("G18").Value = 0.00271495312305698
$ws.Range("I18").Value = 3.988361548184301
$ws.Range("J18").Value = 0.01688891990283281
$ws.Range("L18").Value = 1.067519217961632
$ws.Range("N18").Value = 3.511436773817081
$ws.Range("B19").Value = 6.683829022027453
$ws.Range("C19").Value = 1.254095577088435
$ws.Range("D19").Value = 0.3424918374226138
$ws.Range("E19").Value = 0.05761442686884477
$ws.Range("G19").Value = 0.002715591082165458
$ws.Range("I19").Value = 3.985664256401137
$ws.Range("J19").Value = 0.01675573453523072
$ws.Range("L19").Value = 1.066473360521314
$ws.Range("N19").Value = 3.511543919689871
$ws.Range("B20").Value = 6.742862144270703
$ws.Range("C20").Value = 1.273036637994096
$ws.Range("D20").Value = 0.3427510237091411
$ws.Range("E20").Value = 0.05755089011083703
$ws.Range("G20").Value = 0.002712736950961304
$ws.Range("I20").Value = 3.997923116202713
$ws.Range("J20").Value = 0.01735505000183224
$ws.Range("L20").Value = 1.071218316034617
$ws.Range("N20").Value = 3.511118738151509
$ws.Range("B21").Value = 6.946115574317673
$ws.Range("C21").Value = 1.337680414592626
$ws.Range("D21").Value = 0.3437899110872138
$ws.Range("E21").Value = 0.05734694639046722
$ws.Range("G21").Value = 0.002703439633140118
$ws.Range("I21").Value = 4.041235419024844
$ws.Range("J21").Value = 0.01936839683837377
$ws.Range("L21").Value = 1.087836898985756
$ws.Range("N21").Value = 3.51069304599929
$ws.Range("B22").Value = 7.082646891438685
$ws.Range("C22").Value = 1.380690008891747
$ws.Range("D22").Value = 0.3445944297800736
$ws.Range("E22").Value = 0.05722073763406232
$ws.Range("G22").Value = 0.002697579776514171
$ws.Range("I22").Value = 4.07113940668097
$ws.Range("J22").Value = 0.02068441954524047
$ws.Range("L22").Value = 1.099205489805172
$ws.Range("N22").Value = 3.511169116854575
$ws.Range("B23").Value = 7.009437403128231
$ws.Range("C23").Value = 1.35766451778494
$ws.Range("D23").Value = 0.3441535363263313
$ws.Range("E23").Value = 0.05728745126772372
$ws.Range("G23").Value = 0.002700687567035445
$ws.Range("I23").Value = 4.055032079398842
$ws.Range("J23").Value = 0.0199819915782129
$ws.Range("L23").Value = 1.093091209869613
$ws.Range("N23").Value = 3.510845568128246
$ws.Range("B24").Value = 6.739601939574641
$ws.Range("C24").Value = 1.271992738397216
$ws.Range("D24").Value = 0.3427361589194362
$ws.Range("E24").Value = 0.05755434332979092
$ws.Range("G24").Value = 0.002712892584052288
$ws.Range("I24").Value = 3.997241961266269
$ws.Range("J24").Value = 0.01732214045414793
$ws.Range("L24").Value = 1.070955210492485
$ws.Range("N24").Value = 3.511138331559181
$ws.Range("B25").Value = 6.464180148476089
$ws.Range("C25").Value = 1.18276871857131
$ws.Range("D25").Value = 0.3417444475173426
$ws.Range("E25").Value = 0.05787278104055704
$ws.Range("G25").Value = 0.002726995433535087
$ws.Range("I25").Value = 3.941673770062323
$ws.Range("J25").Value = 0.01445135385478125
$ws.Range("L25").Value = 1.049234337598989
$ws.Range("N25").Value = 3.514648474105002
